$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates as Excel serial numbers, matching existing column A format)
$data = @(
    @(43812, 0, 0, 0, 0, 12.5, 0, 0, 0, 0, 0, 0, 3),
    @(43813, 0, 5, 2.5, 5, 0, 3.5, 25, 0, 0, 1, 2, 2),
    @(43814, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 2, 2)
)

$startRow = 24

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Reuse column A's existing date number formatting (style) for the new rows
# by copying the format from the row directly above each new row.
$ws.Range("A23").Copy()
$ws.Range("A24:A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to match the new last populated cell
$ws.Range("K26").Select()
